# time measurement adj meth
# Remove the "Testdatei in Main (...)" and "code aufraeumen und kommentieren"
# to-do entries, and move the active cell selection down to B19.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the two removed to-do entries (their rows disappear from sheetData
# entirely once they have no remaining cell content).
$ws.Range("B7").ClearContents()
$ws.Range("B11").ClearContents()

# Update the saved selection/active cell for the sheet view.
$ws.Range("B19").Select() | Out-Null
